$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns for each coin row to reflect
# the latest scrape from the GitHub Actions cron job.

$ws.Range("D2").Value = "26.167.78"
$ws.Range("E2").Value = "  +0.06%  "

$ws.Range("D3").Value = "1.650.99"
$ws.Range("E3").Value = "  -0.36%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.014"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.37%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.72%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5061"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.73%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.014"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.41%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2591"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.18%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06469"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.77%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.60"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.26%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07773"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.04%  "

$ws.Range("D12").Value = "1.675.84"
$ws.Range("E12").Value = "  +1.08%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.284"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.31%  "

$ws.Range("D14").Value = "1.871.16"
$ws.Range("E14").Value = "  -0.74%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5495"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.60%  "

$ws.Range("D16").Value = "0.0₅7973"
$ws.Range("E16").Value = "  -0.45%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.00"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.47%  "

$ws.Range("D18").Value = "26.122.88"
$ws.Range("E18").Value = "  -0.23%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.013"
$ws.Range("D19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "204.74"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.45%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.337"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.25%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.05"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.21%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.003"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.43%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.014"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.37%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.965"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +11.81%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.58"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.22%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1163"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.44%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.84%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.773"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.48%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05093"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.92%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.249"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.80%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.286"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.19%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.215"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.13%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.555"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.46%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.355"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.11%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9026"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.22%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.630"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -4.15%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5678"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.42%  "

$ws.Range("D39").Value = "1.154.83"
$ws.Range("E39").Value = "  +0.19%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01580"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.63%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.584"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.09%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.015"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.39%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.690"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.63%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8218"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.89%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "100.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.32%  "

$ws.Range("D46").Value = "1.782.41"
$ws.Range("E46").Value = "  -0.85%  "

$ws.Range("D47").Value = "0.0₈114"
$ws.Range("E47").Value = "  +4.24%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4556"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.05%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.012"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.23%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "55.27"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.07%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05055"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.64%  "
